$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - column headers (nicer, capitalized display names)
$ws.Range("B1").Value = "Incendios"
$ws.Range("C1").Value = "Comarca nombre"
$ws.Range("D1").Value = "Comarca código"
$ws.Range("E1").Value = "Superficie forestal afectada"
$ws.Range("F1").Value = "CCAA código"
$ws.Range("G1").Value = "Provincia código"
$ws.Range("H1").Value = "Municipio código"
$ws.Range("I1").Value = "Provincia nombre"
$ws.Range("J1").Value = "Año"
$ws.Range("K1").Value = "Municipio nombre"

# Row 2 - sdmx dimension / measure annotations
$ws.Range("A2").Value = "iaest-measure:"
$ws.Range("B2").Value = "iaest-measure:incendios"
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("D2").Value = "null"
$ws.Range("E2").Value = "iaest-measure:superficie-forestal-afectada"
$ws.Range("F2").Value = "null"
$ws.Range("G2").Value = "null"
$ws.Range("H2").Value = "null"
$ws.Range("I2").Value = "sdmx-dimension:refArea"
$ws.Range("J2").Value = "sdmx-dimension:refPeriod"
$ws.Range("K2").Value = "sdmx-dimension:refArea"

# Row 3 - medida / dim markers
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "null"
$ws.Range("G3").Value = "null"
$ws.Range("H3").Value = "null"
$ws.Range("I3").Value = "dim"
$ws.Range("J3").Value = "dim"
$ws.Range("K3").Value = "dim"

# Row 4 - data types / codelist URIs
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "URI-comarca"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "xsd:double"
$ws.Range("F4").Value = "null"
$ws.Range("G4").Value = "null"
$ws.Range("H4").Value = "null"
$ws.Range("I4").Value = "URI-Provincia"
$ws.Range("J4").Value = "xsd:date"
$ws.Range("K4").Value = "URI-Municipio"

# Row 5 - mapping file annotation moved from C5 to J5
$ws.Range("C4").Copy()
$ws.Range("J5").PasteSpecial(-4122)
$ws.Range("J5").Value = "mapping-ano.xlsx"
$ws.Range("C5").Clear()
